# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values on the zh-cn and
# de-de report sheets for the 34e60af9-... row (row 4, also mirrored
# on row 5) to reflect the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-11 14:14:38"
$zhcn.Range("E5").Value = "2016-03-11 14:14:38"
$zhcn.Range("H4").Value = "2016-03-11 14:14:59"
$zhcn.Range("H5").Value = "2016-03-11 14:14:59"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-11 14:14:42"
$dede.Range("E5").Value = "2016-03-11 14:14:42"
$dede.Range("H4").Value = "2016-03-11 14:15:08"
$dede.Range("H5").Value = "2016-03-11 14:15:08"
